$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("progress")

# --- 1. Insert a new row at position 3 (shifts old rows 3..24 down to 4..25) ---
$ws.Rows.Item(3).Insert()

# Fix A2's formula: Excel's auto-adjust during insert pushed the 1048576 boundary to
# 1048577 (out of range); restore the intended range explicitly.
$ws.Range("A2").Formula = "=AVERAGE(A5:A1048576)"

# --- 2. Remove the old "Advanced ssd1306 functions" snippets from column E ---
# (old rows 4-8 are now rows 5-9; keep their styles, just clear the text)
$ws.Range("E5:E9").ClearContents()

# --- 3. Clear the old column G/N helper block; it is being rebuilt in columns F/G ---
$ws.Range("G11:G25").ClearContents()
$ws.Range("N11:N25").ClearContents()

# --- 4. New header block (rows 3-4, columns F/G) ---
$ws.Range("F3").Value2 = "/*"
$ws.Range("F4").Value2 = "Different sub-projects to port the Arduino project to ESP-IDF"
$ws.Range("G4").Value2 = "Different projects were created to port the Arduino version to ESP-IDF:"

# --- 5. Update the text for what is now row 16 (Bedclock_IDF_V12 description) ---
$ws.Range("B16").Value2 = "First working version with graphics layout"

# --- 6. Rebuild the two helper columns as shared formulas F5:F16 / G5:G16 ---
$ws.Range("F5:F16").Formula = '="    "&C5&" : "&B5'
$ws.Range("G5:G16").Formula = '="  * `"&C5&"` : "&B5'

# --- 7. Closing comment marker, now in column F ---
$ws.Range("F17").Value2 = "*/"

# --- 8. Restore selection / active cell like the saved workbook ---
$ws.Range("G4").Select()

$wb.Save()
